# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets to the latest scraped values.
# Row 26 differs slightly between the two sheets in the source data
# (823 vs 824), matching the generated export.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

$sheetExhibition.Range("F2").Value = 861
$sheetExhibition.Range("F3").Value = 612
$sheetExhibition.Range("F4").Value = 2171
$sheetExhibition.Range("F6").Value = 12652
$sheetExhibition.Range("F13").Value = 13667
$sheetExhibition.Range("F14").Value = 14012
$sheetExhibition.Range("F23").Value = 1054
$sheetExhibition.Range("F26").Value = 823
$sheetExhibition.Range("F27").Value = 5134
$sheetExhibition.Range("F28").Value = 6
$sheetExhibition.Range("F29").Value = 256

$sheetAllTypes.Range("F2").Value = 861
$sheetAllTypes.Range("F3").Value = 612
$sheetAllTypes.Range("F4").Value = 2171
$sheetAllTypes.Range("F6").Value = 12652
$sheetAllTypes.Range("F13").Value = 13667
$sheetAllTypes.Range("F14").Value = 14012
$sheetAllTypes.Range("F23").Value = 1054
$sheetAllTypes.Range("F26").Value = 824
$sheetAllTypes.Range("F27").Value = 5134
$sheetAllTypes.Range("F28").Value = 6
$sheetAllTypes.Range("F29").Value = 256
